$d = $word.ActiveDocument

# --- Header fields: re-assert the existing text so Word recombines the
#     runs that were split for spell-check tracking and drops the now
#     stale w:proofErr markers around them. ---
$d.Content.Find.Execute("COURSE: [data.course_name]", $true, $false, $false, $false, $false, $true, 1, $false, "COURSE: [data.course_name]", 2)
$d.Content.Find.Execute("UNIT: [data.unit_name]", $true, $false, $false, $false, $false, $true, 1, $false, "UNIT: [data.unit_name]", 2)
$d.Content.Find.Execute("SESSION: [data.session], [data.year]", $true, $false, $false, $false, $false, $true, 1, $false, "SESSION: [data.session], [data.year]", 2)
$d.Content.Find.Execute("STD UNITS: [data.std_units]", $true, $false, $false, $false, $false, $true, 1, $false, "STD UNITS: [data.std_units]", 2)
$d.Content.Find.Execute("CLASS CODE: [data.class_code]", $true, $false, $false, $false, $false, $true, 1, $false, "CLASS CODE: [data.class_code]", 2)
$d.Content.Find.Execute("TEACHER: [data.teacher_name]", $true, $false, $false, $false, $false, $true, 1, $false, "TEACHER: [data.teacher_name]", 2)

# --- Unit goals / content merge placeholders: now resolved against the
#     ".text" property of the block item rather than the item itself. ---
$d.Content.Find.Execute("[unit_goals; block=tbs:listitem]", $true, $false, $false, $false, $false, $true, 1, $false, "[unit_goals.text; block=tbs:listitem]", 2)
$d.Content.Find.Execute("[unit_content;block=tbs:listitem]", $true, $false, $false, $false, $false, $true, 1, $false, "[unit_content.text;block=tbs:listitem]", 2)

# --- Assessment table cells: re-assert text to drop stale proofErr runs. ---
$d.Content.Find.Execute("[assessment.name; block=tbs:row]", $true, $false, $false, $false, $false, $true, 1, $false, "[assessment.name; block=tbs:row]", 2)
$d.Content.Find.Execute("[assessment.weighting]", $true, $false, $false, $false, $false, $true, 1, $false, "[assessment.weighting]", 2)
$d.Content.Find.Execute("[assessment.due_date]", $true, $false, $false, $false, $false, $true, 1, $false, "[assessment.due_date]", 2)

# --- Move the "_GoBack" bookmark (Word's "last edit position" marker) to
#     sit right after ".text" in the CONTENT placeholder, matching where
#     the author's final edit in this session landed. ---
$r = $d.Content
$r.Find.Execute("unit_content.text", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pos = $r.End
$bmRange = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $bmRange)
